# Add a new "Hungary" worksheet (cloned from "Slovakia") with its own
# market data, matching the author's "Added HungaryFC Test data" commit.

$wb = $excel.ActiveWorkbook

# --- 1. Clone the Slovakia sheet and place it after the last sheet -------
$template = $wb.Worksheets.Item("Slovakia")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy([System.Reflection.Missing]::Value, $lastSheet)

$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"

# --- 2. Fill in the Hungary-specific values -------------------------------
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-4308/T3593/T3618"

# Match the thin-border style already used on B3/B4 elsewhere in the sheet.
$hungary.Range("B3").Copy()
$hungary.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Widen column B so the longer "Hungary Market" text fits -----------
# (mirrors Excel's "best fit" column-width behaviour for the new text)
$hungary.Columns.Item(2).ColumnWidth = 20.25

# --- 4. View state: Hungary tab becomes the active/selected tab -----------
$hungary.Range("B7").Select()

# The UK sheet's selection becomes a "select all" (A1:XFD1048576),
# and it is no longer the selected tab (Hungary is now).
$uk = $wb.Worksheets.Item("UK")
$uk.Cells.Select()

# Make sure Hungary ends up as the active sheet/tab.
$hungary.Activate()
$hungary.Range("B7").Select()
